# Update the LR-pair TPM results table (rows 2-7, columns A:T) with the
# refreshed per-(sending cluster, target cluster) statistics computed from
# the new TPM input. Row data below mirrors the updated dataset row-by-row;
# new cluster-pair combinations are appended as rows 6 and 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Ccl5", "Gpr75", "ECs", 1, 0.3333333333333333, 0.021814, 0.065442, 0.1008129179549036, 0.1008129179549036, 1, 0.3333333333333333, 0.03650533333333333, 0.109516, 0.06499111027765844, 0.06499111027765843, 0.0007963273413333333, 0.007166946072, 0.006551943468219674, 0.006551943468219672),
    @("ECs", "Ccl5", "Gpr75", "FAPs", 1, 0.3333333333333333, 0.021814, 0.065442, 0.1008129179549036, 0.1008129179549036, 2, 0.6666666666666666, 0.2483593333333333, 0.745078, 0.4421586477177508, 0.4421586477177507, 0.005417710497333333, 0.048759394476, 0.04457530347542074, 0.04457530347542073),
    @("ECs", "Ccl5", "Gpr75", "MuSCs", 1, 0.3333333333333333, 0.021814, 0.065442, 0.1008129179549036, 0.1008129179549036, 2, 0.6666666666666666, 0.2768326666666667, 0.830498, 0.4928502420045909, 0.4928502420045908, 0.006038827790666667, 0.05434945011599999, 0.04968567101126322, 0.04968567101126321),
    @("MuSCs", "Ccl5", "Gpr75", "ECs", 1, 0.3333333333333333, 0.194567, 0.583701, 0.8991870820450963, 0.8991870820450963, 1, 0.3333333333333333, 0.03650533333333333, 0.109516, 0.06499111027765844, 0.06499111027765843, 0.007102733190666667, 0.06392459871600001, 0.05843916680943877, 0.05843916680943875),
    @("MuSCs", "Ccl5", "Gpr75", "FAPs", 1, 0.3333333333333333, 0.194567, 0.583701, 0.8991870820450963, 0.8991870820450963, 2, 0.6666666666666666, 0.2483593333333333, 0.745078, 0.4421586477177508, 0.4421586477177507, 0.04832253040866667, 0.434902773678, 0.3975833442423301, 0.3975833442423299),
    @("MuSCs", "Ccl5", "Gpr75", "MuSCs", 1, 0.3333333333333333, 0.194567, 0.583701, 0.8991870820450963, 0.8991870820450963, 2, 0.6666666666666666, 0.2768326666666667, 0.830498, 0.4928502420045909, 0.4928502420045908, 0.05386250145533334, 0.484762513098, 0.4431645709933277, 0.4431645709933276)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $rows[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$col - 1]
    }
}

